# The two data records currently stored in rows 2 and 3 of "Artfynd" get
# swapped: row 2 ends up holding what used to be row 3's record, and row 3
# ends up holding what used to be row 2's record. Columns C, D, P, S, T, U,
# V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY already hold identical
# values in both rows, so only the columns below actually change value:
# A, B, E, F, G, H, I, Q, R (plus the empty "Bestamningsmetod" (AF) marker
# cell, which is present on whichever row currently has no "Antal" (I)
# figure).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Read with Value2 - this interop shim's parameterless `.Value` getter
# does not marshal back correctly when read (assigning to `.Value` works
# fine, it is only the read side that is unreliable), so `.Value2` is used
# for every read below.

# --- Capture row 2 (old) values ---
$A2 = $ws.Cells.Item(2, 1).Value2
$B2 = $ws.Cells.Item(2, 2).Value2
$E2 = $ws.Cells.Item(2, 5).Value2
$F2 = $ws.Cells.Item(2, 6).Value2
$G2 = $ws.Cells.Item(2, 7).Value2
$H2 = $ws.Cells.Item(2, 8).Value2
$I2 = $ws.Cells.Item(2, 9).Value2
$Q2 = $ws.Cells.Item(2, 17).Value2
$R2 = $ws.Cells.Item(2, 18).Value2

# --- Capture row 3 (old) values ---
$A3 = $ws.Cells.Item(3, 1).Value2
$B3 = $ws.Cells.Item(3, 2).Value2
$E3 = $ws.Cells.Item(3, 5).Value2
$F3 = $ws.Cells.Item(3, 6).Value2
$G3 = $ws.Cells.Item(3, 7).Value2
$H3 = $ws.Cells.Item(3, 8).Value2
$I3 = $ws.Cells.Item(3, 9).Value2
$Q3 = $ws.Cells.Item(3, 17).Value2
$R3 = $ws.Cells.Item(3, 18).Value2

# --- Write old row 3's data into row 2 ---
$ws.Cells.Item(2, 1).Value = $A3
$ws.Cells.Item(2, 2).Value = $B3
$ws.Cells.Item(2, 5).Value = $E3
$ws.Cells.Item(2, 6).Value = $F3
$ws.Cells.Item(2, 7).Value = $G3
$ws.Cells.Item(2, 8).Value = $H3
# "Antal" (I) keeps being stored as text, same as it was on the row it
# came from - force text so the numeric-looking "30" isn't auto-converted
# to a number.
$ws.Cells.Item(2, 9).NumberFormat = "@"
$ws.Cells.Item(2, 9).Value = $I3
$ws.Cells.Item(2, 17).Value = $Q3
$ws.Cells.Item(2, 18).Value = $R3

# --- Write old row 2's data into row 3 ---
$ws.Cells.Item(3, 1).Value = $A2
$ws.Cells.Item(3, 2).Value = $B2
$ws.Cells.Item(3, 5).Value = $E2
$ws.Cells.Item(3, 6).Value = $F2
$ws.Cells.Item(3, 7).Value = $G2
$ws.Cells.Item(3, 8).Value = $H2
$ws.Cells.Item(3, 9).NumberFormat = "@"
$ws.Cells.Item(3, 9).Value = $I2
$ws.Cells.Item(3, 17).Value = $Q2
$ws.Cells.Item(3, 18).Value = $R2

# --- The empty "Bestamningsmetod" (AF) marker follows the empty "Antal"
# (I) cell: it sits on row 2 when row 2 has no Antal figure, and on row 3
# when row 3 has no Antal figure. Before the edit row 2 had the Antal
# figure (so AF2 was the empty marker); after swapping, row 3 has no
# Antal figure, so the empty marker moves to AF3.
$ws.Cells.Item(2, 32).Value = ""
$ws.Cells.Item(3, 32).Value = ""
